$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update login values (A2: AdminID value, C2: LoginLabel value)
$ws.Range("A2").Value = "testSean"
$ws.Range("C2").Value = "SEANTEST PROCTOR"

# Update the active selection to C2
$ws.Range("C2").Select()
